$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Mining)
$ws.Range("B2").Value = "0.006`n (0.020)"
$ws.Range("C2").Value = "0.043`n (0.053)"
$ws.Range("D2").Value = "0.010`n (0.014)"
$ws.Range("E2").Value = "0.002`n (0.034)"
$ws.Range("F2").Value = "-0.010`n (0.011)"
$ws.Range("G2").Value = "0.013`n (0.024)"

# Row 3 (Partisanship)
$ws.Range("B3").Value = "-0.366***`n (0.115)"
$ws.Range("C3").Value = "-0.594**`n (0.226)"
$ws.Range("D3").Value = "0.472***`n (0.082)"
$ws.Range("E3").Value = "0.240`n (0.147)"
$ws.Range("F3").Value = "0.341***`n (0.062)"
$ws.Range("G3").Value = "0.137`n (0.100)"

# Row 4 (Deregulated)
$ws.Range("B4").Value = "0.311*`n (0.181)"
$ws.Range("D4").Value = "0.079`n (0.126)"
$ws.Range("F4").Value = "0.138`n (0.102)"

# Row 5 (R-squared)
$ws.Range("B5").Value = 0.2839416083222908
$ws.Range("C5").Value = 0.1533856611860601
$ws.Range("D5").Value = 0.4885907484925434
$ws.Range("E5").Value = 0.06683910056541831
$ws.Range("F5").Value = 0.3711881025611223
$ws.Range("G5").Value = 0.04549049054619092

# Row 6 (N)
$ws.Range("B6").Value = 56
$ws.Range("C6").Value = 56
$ws.Range("D6").Value = 49
$ws.Range("E6").Value = 49
$ws.Range("F6").Value = 57
$ws.Range("G6").Value = 57
